# Apply cryptos-list refresh (prices + 1h volume deltas; two coin rows re-sorted)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.973.27"
$ws.Range("E2").Value = "  -3.67%  "

# Row 3
$ws.Range("D3").Value = "3.343.23"
$ws.Range("E3").Value = "  -0.80%  "

# Row 5
$ws.Range("D5").Value = "'574.90"
$ws.Range("E5").Value = "  -3.06%  "

# Row 6
$ws.Range("D6").Value = "183.03"
$ws.Range("E6").Value = "  -4.99%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  -1.52%  "

# Row 9
$ws.Range("E9").Value = "  -3.52%  "

# Row 10
$ws.Range("E10").Value = "  -1.56%  "

# Row 11
$ws.Range("D11").Value = "0.404"
$ws.Range("E11").Value = "  -4.32%  "

# Row 12
$ws.Range("D12").Value = "3.923.77"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13
$ws.Range("E13").Value = "  -0.90%  "

# Row 14
$ws.Range("D14").Value = "27.26"
$ws.Range("E14").Value = "  -4.90%  "

# Row 15
$ws.Range("D15").Value = "67.021.86"
$ws.Range("E15").Value = "  -3.62%  "

# Row 16
$ws.Range("E16").Value = "  -2.45%  "

# Row 17
$ws.Range("D17").Value = "3.349.76"
$ws.Range("E17").Value = "  -0.53%  "

# Row 18
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "437.21"
$ws.Range("E18").Value = "  -2.76%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "13.72"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20
$ws.Range("D20").Value = "5.71"
$ws.Range("E20").Value = "  -2.38%  "

# Row 21
$ws.Range("D21").Value = "7.66"
$ws.Range("E21").Value = "  -2.16%  "

# Row 22
$ws.Range("D22").Value = "73.82"
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").Value = "0.524"
$ws.Range("E24").Value = "  +0.96%  "

# Row 25
$ws.Range("E25").Value = "  -2.71%  "

# Row 26
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  -4.26%  "

# Row 28
$ws.Range("E28").Value = "  -0.14%  "

# Row 30
$ws.Range("D30").Value = "22.89"
$ws.Range("E30").Value = "  -1.64%  "

# Row 31
$ws.Range("D31").Value = "5.36"
$ws.Range("E31").Value = "  -4.58%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.25"
$ws.Range("E32").Value = "  -3.80%  "

# Row 33
$ws.Range("D33").Value = "6.85"
$ws.Range("E33").Value = "  -2.67%  "

# Row 34
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.05%  "

# Row 35
$ws.Range("E35").Value = "  -1.19%  "

# Row 36
$ws.Range("D36").Value = "160.52"
$ws.Range("E36").Value = "  -2.75%  "

# Row 37
$ws.Range("D37").Value = "27.59"
$ws.Range("E37").Value = "  +1.09%  "

# Row 38
$ws.Range("E38").Value = "  -4.97%  "

# Row 39
$ws.Range("D39").Value = "2.837.56"
$ws.Range("E39").Value = "  +3.46%  "

# Row 40
$ws.Range("D40").Value = "0.794"
$ws.Range("E40").Value = "  -3.20%  "

# Row 41
$ws.Range("D41").Value = "4.46"
$ws.Range("E41").Value = "  -3.30%  "

# Row 42
$ws.Range("E42").Value = "  -4.58%  "

# Row 43
$ws.Range("D43").Value = "0.0678"
$ws.Range("E43").Value = "  -1.97%  "

# Row 44
$ws.Range("D44").Value = "40.32"
$ws.Range("E44").Value = "  -1.23%  "

# Row 45
$ws.Range("D45").Value = "24.75"
$ws.Range("E45").Value = "  -3.57%  "

# Row 46
$ws.Range("D46").Value = "2.38"
$ws.Range("E46").Value = "  -6.52%  "

# Row 47
$ws.Range("D47").Value = "324.55"
$ws.Range("E47").Value = "  -5.58%  "

# Row 48
$ws.Range("E48").Value = "  -4.06%  "

# Row 49
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "31.96"
$ws.Range("E49").Value = "  -3.25%  "

# Row 50
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "0.992"
$ws.Range("E50").Value = "  -4.26%  "

# Row 51
$ws.Range("D51").Value = "6.17"
$ws.Range("E51").Value = "  -2.61%  "
